$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1072.1666
$ws.Range("I28").Value = 328.5
$ws.Range("J28").Value = 2559.5
$ws.Range("K28").Value = 328.5
$ws.Range("L28").Value = 2559.5
$ws.Range("M28").Value = 156.5
$ws.Range("N28").Value = -3529.5

$ws.Range("H129").Value = 1651.2106
$ws.Range("I129").Value = 514.2222
$ws.Range("K129").Value = 1542.6666
$ws.Range("M129").Value = 3457.3334

$ws.Range("H132").Value = 4448.6177
$ws.Range("I132").Value = 4727.516
$ws.Range("K132").Value = 14182.548
$ws.Range("M132").Value = -11652.548

$ws.Range("H135").Value = 733.5625
$ws.Range("I135").Value = 371.57144
$ws.Range("K135").Value = 3344.14296
$ws.Range("M135").Value = -809.1429600000001

$ws.Range("H137").Value = 5407.3687
$ws.Range("I137").Value = 2346
$ws.Range("J137").Value = 15271.777
$ws.Range("K137").Value = 7038
$ws.Range("L137").Value = 45815.331
$ws.Range("M137").Value = -4488
$ws.Range("N137").Value = -50915.331

$ws.Range("H138").Value = 2086.0588
$ws.Range("I138").Value = 996.7273
$ws.Range("J138").Value = 4083.1667
$ws.Range("K138").Value = 2990.1819
$ws.Range("L138").Value = 12249.5001
$ws.Range("M138").Value = 2149.8181
$ws.Range("N138").Value = -22529.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1171.8518
$ws.Range("I2").Value = 1172.7333
$ws.Range("J2").Value = 1170.75
$ws.Range("K2").Value = 1172.7333
$ws.Range("L2").Value = 1170.75
$ws.Range("M2").Value = -1059.7333
$ws.Range("N2").Value = -1396.75

$ws.Range("H24").Value = 36142.43
$ws.Range("J24").Value = 36142.43
$ws.Range("L24").Value = 36142.43
$ws.Range("N24").Value = -36890.43

$ws.Range("H52").Value = 49999
$ws.Range("J52").Value = 49999
$ws.Range("L52").Value = 49999
$ws.Range("N52").Value = -50635

$ws.Range("H61").Value = 2353.658
$ws.Range("I61").Value = 1613.125
$ws.Range("K61").Value = 1613.125
$ws.Range("M61").Value = -1401.125

$ws.Range("H74").Value = 152037.08
$ws.Range("I74").Value = 232865.2
$ws.Range("J74").Value = 2815.923
$ws.Range("K74").Value = 232865.2
$ws.Range("L74").Value = 2815.923
$ws.Range("M74").Value = -231991.2
$ws.Range("N74").Value = -4563.923

$ws.Range("H77").Value = 152037.08
$ws.Range("I77").Value = 232865.2
$ws.Range("J77").Value = 2815.923
$ws.Range("K77").Value = 1164326
$ws.Range("L77").Value = 14079.615
$ws.Range("M77").Value = -1159958
$ws.Range("N77").Value = -22815.615

$ws.Range("H100").Value = 36142.43
$ws.Range("J100").Value = 36142.43
$ws.Range("L100").Value = 36142.43
$ws.Range("N100").Value = -38306.43

$ws.Range("H116").Value = 1171.8518
$ws.Range("I116").Value = 1172.7333
$ws.Range("J116").Value = 1170.75
$ws.Range("K116").Value = 1172.7333
$ws.Range("L116").Value = 1170.75
$ws.Range("M116").Value = 1121.2667
$ws.Range("N116").Value = -5758.75

$ws.Range("H122").Value = 2804.2856
$ws.Range("I122").Value = 3026.1
$ws.Range("K122").Value = 9078.299999999999
$ws.Range("M122").Value = -6628.299999999999

$ws.Range("H132").Value = 2823.375
$ws.Range("J132").Value = 3376.75
$ws.Range("L132").Value = 10130.25
$ws.Range("N132").Value = -15190.25

$ws.Range("H136").Value = 2353.658
$ws.Range("I136").Value = 1613.125
$ws.Range("K136").Value = 4839.375
$ws.Range("M136").Value = -2289.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1171.8518
$ws.Range("I3").Value = 1172.7333
$ws.Range("J3").Value = 1170.75
$ws.Range("K3").Value = 1172.7333
$ws.Range("L3").Value = 1170.75
$ws.Range("M3").Value = -1058.7333
$ws.Range("N3").Value = -1398.75

$ws.Range("H20").Value = 13892359
$ws.Range("I20").Value = 19234676
$ws.Range("J20").Value = 2333.4
$ws.Range("K20").Value = 19234676
$ws.Range("L20").Value = 2333.4
$ws.Range("M20").Value = -19234429
$ws.Range("N20").Value = -2827.4

$ws.Range("H55").Value = 30780
$ws.Range("J55").Value = 30780
$ws.Range("L55").Value = 30780
$ws.Range("N55").Value = -31326

$ws.Range("H64").Value = 1569.4
$ws.Range("I64").Value = 1399.875
$ws.Range("J64").Value = 2247.5
$ws.Range("K64").Value = 1399.875
$ws.Range("L64").Value = 2247.5
$ws.Range("M64").Value = -1174.875
$ws.Range("N64").Value = -2697.5

$ws.Range("H67").Value = 1569.4
$ws.Range("I67").Value = 1399.875
$ws.Range("J67").Value = 2247.5
$ws.Range("K67").Value = 1399.875
$ws.Range("L67").Value = 2247.5
$ws.Range("M67").Value = -619.875
$ws.Range("N67").Value = -3807.5

$ws.Range("H68").Value = 15000
$ws.Range("J68").Value = 15000
$ws.Range("L68").Value = 15000
$ws.Range("N68").Value = -16622

$ws.Range("H71").Value = 15000
$ws.Range("J71").Value = 15000
$ws.Range("L71").Value = 45000
$ws.Range("N71").Value = -53112

$ws.Range("H134").Value = 2341.1365
$ws.Range("I134").Value = 1968.9333
$ws.Range("K134").Value = 5906.7999
$ws.Range("M134").Value = -3371.7999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2781885.8
$ws.Range("I31").Value = 3105.875
$ws.Range("K31").Value = 3105.875
$ws.Range("M31").Value = -2810.875

$ws.Range("H34").Value = 2781885.8
$ws.Range("I34").Value = 3105.875
$ws.Range("K34").Value = 3105.875
$ws.Range("M34").Value = -2903.875

$ws.Range("H58").Value = 2649.4546
$ws.Range("I58").Value = 2500
$ws.Range("J58").Value = 2682.6667
$ws.Range("K58").Value = 2500
$ws.Range("L58").Value = 2682.6667
$ws.Range("M58").Value = -2297
$ws.Range("N58").Value = -3088.6667

$ws.Range("H136").Value = 2649.4546
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 2682.6667
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 8048.000100000001
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -13148.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 572.8889
$ws.Range("I14").Value = 572.8889
$ws.Range("K14").Value = 1718.6667
$ws.Range("M14").Value = -1545.6667

$ws.Range("H41").Value = 1600

$ws.Range("H70").Value = 800
$ws.Range("I70").Value = 800
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2400
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2085
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 800
$ws.Range("I73").Value = 800
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2400
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1308
$ws.Range("N73").ClearContents()

$ws.Range("H122").Value = 893.5
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 18000
$ws.Range("N122").Value = -22900

$ws.Range("H131").Value = 6232.3076
$ws.Range("I131").Value = 10481.5
$ws.Range("J131").Value = 2590.1428
$ws.Range("K131").Value = 31444.5
$ws.Range("L131").Value = 7770.428400000001
$ws.Range("M131").Value = -26404.5
$ws.Range("N131").Value = -17850.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13249.25
$ws.Range("I43").Value = 3999
$ws.Range("J43").Value = 41000
$ws.Range("K43").Value = 3999
$ws.Range("L43").Value = 41000
$ws.Range("M43").Value = -3848
$ws.Range("N43").Value = -41302

$ws.Range("H70").Value = 4807.2144
$ws.Range("I70").Value = 2866.3333
$ws.Range("K70").Value = 2866.3333
$ws.Range("M70").Value = -2596.3333

$ws.Range("H73").Value = 4807.2144
$ws.Range("I73").Value = 2866.3333
$ws.Range("K73").Value = 2866.3333
$ws.Range("M73").Value = -1930.3333

$ws.Range("H126").Value = 8812.857
$ws.Range("J126").Value = 11436.8
$ws.Range("L126").Value = 34310.39999999999
$ws.Range("N126").Value = -39250.39999999999

$ws.Range("H132").Value = 2048.1667
$ws.Range("I132").Value = 1571.95
$ws.Range("K132").Value = 4715.85
$ws.Range("M132").Value = -2185.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 25625
$ws.Range("I40").Value = 48250
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 48250
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -48114
$ws.Range("N40").Value = -3272

$ws.Range("H46").Value = 2942.4707
$ws.Range("I46").Value = 2248.889
$ws.Range("J46").Value = 3722.75
$ws.Range("K46").Value = 2248.889
$ws.Range("L46").Value = 3722.75
$ws.Range("M46").Value = -2060.889
$ws.Range("N46").Value = -4098.75

$ws.Range("H107").Value = 4998
$ws.Range("I107").Value = 4998
$ws.Range("K107").Value = 4998
$ws.Range("M107").Value = -3078

$ws.Range("H132").Value = 5247.6895
$ws.Range("I132").Value = 4618.3335
$ws.Range("J132").Value = 5691.9414
$ws.Range("K132").Value = 13855.0005
$ws.Range("L132").Value = 17075.8242
$ws.Range("M132").Value = -11325.0005
$ws.Range("N132").Value = -22135.8242

$ws.Range("H133").Value = 87249.5
$ws.Range("J133").Value = 87249.5
$ws.Range("L133").Value = 87249.5
$ws.Range("N133").Value = -92309.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1312.2
$ws.Range("I96").Value = 1312.2
$ws.Range("K96").Value = 1312.2
$ws.Range("M96").Value = 60.79999999999995

$ws.Range("H126").Value = 1178.8
$ws.Range("I126").Value = 1096.3334
$ws.Range("J126").Value = 1302.5
$ws.Range("K126").Value = 3289.0002
$ws.Range("L126").Value = 3907.5
$ws.Range("M126").Value = -819.0001999999999
$ws.Range("N126").Value = -8847.5

$ws.Range("H136").Value = 4180.2334
$ws.Range("I136").Value = 2213.037
$ws.Range("J136").Value = 21885
$ws.Range("K136").Value = 6639.110999999999
$ws.Range("L136").Value = 65655
$ws.Range("M136").Value = -4089.110999999999
$ws.Range("N136").Value = -70755
